# Import Job Layout: show Job Available to Import and Duplicated Job
# Adds three new job rows to the bottom of the Jobs list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newJobs = @(
    @{ Num = "J99-9999"; Name = "TEST EXCEL" },
    @{ Num = "J99-9998"; Name = "TEST EXCEL 2" },
    @{ Num = "J99-9997"; Name = "TEST EXCEL 3" }
)

$startRow = 119

for ($i = 0; $i -lt $newJobs.Count; $i++) {
    $row = $startRow + $i
    $job = $newJobs[$i]

    $ws.Cells.Item($row, 1).Value = $job.Num
    $ws.Cells.Item($row, 2).Value = $job.Name
    $ws.Cells.Item($row, 2).HorizontalAlignment = -4131
    $ws.Cells.Item($row, 3).Value = 2021
}

$ws.Range("B121").Select()
